# [IMP] Update import guest
# Swap the "English Name" and "Khmer Name" columns (B and C) — headers and data —
# and apply a couple of small formatting touch-ups that came along with the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap column B and column C content for rows 1, 3 and 4 (row 2 is blank) ---
$rows = @(1, 3, 4)
foreach ($r in $rows) {
    $bCell = $ws.Cells.Item($r, 2)
    $cCell = $ws.Cells.Item($r, 3)
    $tmp = $bCell.Value()
    $bCell.Value = $cCell.Value()
    $cCell.Value = $tmp
}

# --- Center-align A4 (new style picked up for the guest-group cell) ---
$ws.Range("A4").HorizontalAlignment = -4108  # xlCenter

# --- New row 6 / cell C6 gets the same look as the row-3 data cells ---
$ws.Range("C3").Copy()
$ws.Range("C6").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Rows.Item(6).RowHeight = 15.75

# --- Move the active selection to C5 (was D5) ---
$ws.Range("C5").Select()
